$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172335505485535
$ws.Range("B1").Value = 2.223217964172363
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.372478008270264
$ws.Range("E1").Value = 1.228035926818848
